$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 44-45, pushing the existing
# rows 44-58 down to 46-60 (preserving their data/styles).
$ws.Range("A44:A45").EntireRow.Insert()

# Populate the first new row (44) with the newest weekly record.
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44529
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112026
$ws.Range("G44").Value = "Haba"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 125
$ws.Range("K44").Value = 8000
$ws.Range("L44").Value = 8000
$ws.Range("M44").Value = 8000
$ws.Range("N44").Value = "`$/saco 25 kilos"
$ws.Range("O44").Value = "Región de La Araucanía"
$ws.Range("P44").Value = 320
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"

# Populate the second new row (45) with the newest weekly record.
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").Value = 44529
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = 100112026
$ws.Range("G45").Value = "Haba"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 175
$ws.Range("K45").Value = 7000
$ws.Range("L45").Value = 8000
$ws.Range("M45").Value = 7371
$ws.Range("N45").Value = "`$/saco 25 kilos"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 295
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"

$ws.Range("A1").Select()
